$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text formatting (quote-prefix) on numeric-looking Price cells so the
# literal text (with trailing zeros / exact formatting) is preserved, matching
# the original inline-string cell content rather than being converted to a number.
$textCells = @("D5", "D6", "D7", "D8", "D9", "D10", "D11", "D12", "D14", "D16", "D17", "D18", "D21", "D22", "D23", "D24", "D25", "D26", "D27", "D28", "D29", "D30", "D31", "D32", "D33", "D34", "D35", "D36", "D38", "D39", "D40", "D41", "D42", "D43", "D44", "D45", "D46", "D47", "D48", "D49", "D50", "D51")
foreach ($cellRef in $textCells) {
    $ws.Range($cellRef).NumberFormat = "@"
}

# Apply the updated values
$ws.Range('D2').Value = '27.857.10'
$ws.Range('E2').Value = '  -0.93%  '
$ws.Range('D3').Value = '1.906.30'
$ws.Range('E3').Value = '  -0.19%  '
$ws.Range('E4').Value = '  +0.00%  '
$ws.Range('D5').Value = '313.41'
$ws.Range('E5').Value = '  -0.88%  '
$ws.Range('D6').Value = '1.004'
$ws.Range('E6').Value = '  +0.03%  '
$ws.Range('D7').Value = '0.5037'
$ws.Range('E7').Value = '  +4.16%  '
$ws.Range('D8').Value = '0.3817'
$ws.Range('E8').Value = '  -0.08%  '
$ws.Range('D9').Value = '0.07280'
$ws.Range('E9').Value = '  -1.02%  '
$ws.Range('D10').Value = '0.9095'
$ws.Range('E10').Value = '  -2.68%  '
$ws.Range('D11').Value = '20.82'
$ws.Range('E11').Value = '  +0.14%  '
$ws.Range('B12').Value = 'TRON'
$ws.Range('C12').Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range('D12').Value = '0.07654'
$ws.Range('E12').Value = '  -1.71%  '
$ws.Range('B13').Value = 'WrappedEther'
$ws.Range('C13').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D13').Value = '1.926.71'
$ws.Range('E13').Value = '  +0.87%  '
$ws.Range('D14').Value = '5.481'
$ws.Range('E14').Value = '  -0.53%  '
$ws.Range('E15').Value = '  -0.54%  '
$ws.Range('D16').Value = '91.43'
$ws.Range('E16').Value = '  -0.12%  '
$ws.Range('D17').Value = '1.005'
$ws.Range('E17').Value = '  +0.02%  '
$ws.Range('D18').Value = '0.000008718'
$ws.Range('E18').Value = '  -1.30%  '
$ws.Range('D20').Value = '27.868.55'
$ws.Range('E20').Value = '  -1.01%  '
$ws.Range('D21').Value = '14.54'
$ws.Range('E21').Value = '  -2.19%  '
$ws.Range('D22').Value = '5.166'
$ws.Range('E22').Value = '  +0.16%  '
$ws.Range('D23').Value = '10.82'
$ws.Range('E23').Value = '  -0.78%  '
$ws.Range('D24').Value = '154.24'
$ws.Range('E24').Value = '  -1.39%  '
$ws.Range('D25').Value = '1.869'
$ws.Range('E25').Value = '  -2.27%  '
$ws.Range('D26').Value = '2.234'
$ws.Range('E26').Value = '  +5.54%  '
$ws.Range('D27').Value = '18.38'
$ws.Range('E27').Value = '  -0.92%  '
$ws.Range('D28').Value = '115.25'
$ws.Range('E28').Value = '  -0.90%  '
$ws.Range('D29').Value = '4.914'
$ws.Range('E29').Value = '  -0.72%  '
$ws.Range('D30').Value = '0.08987'
$ws.Range('E30').Value = '  +0.70%  '
$ws.Range('D31').Value = '3.211'
$ws.Range('E31').Value = '  -3.86%  '
$ws.Range('D32').Value = '1.234'
$ws.Range('E32').Value = '  -1.56%  '
$ws.Range('D33').Value = '0.7653'
$ws.Range('E33').Value = '  -0.49%  '
$ws.Range('D34').Value = '4.641'
$ws.Range('E34').Value = '  -0.78%  '
$ws.Range('D35').Value = '0.02056'
$ws.Range('E35').Value = '  +0.16%  '
$ws.Range('D36').Value = '2.546'
$ws.Range('E36').Value = '  -2.58%  '
$ws.Range('E37').Value = '  -0.55%  '
$ws.Range('D38').Value = '0.5569'
$ws.Range('E38').Value = '  +1.32%  '
$ws.Range('D39').Value = '3.019'
$ws.Range('E39').Value = '  +1.38%  '
$ws.Range('D40').Value = '0.05257'
$ws.Range('E40').Value = '  -1.02%  '
$ws.Range('D41').Value = '6.972'
$ws.Range('E41').Value = '  -0.56%  '
$ws.Range('D42').Value = '8.486'
$ws.Range('E42').Value = '  +0.27%  '
$ws.Range('D43').Value = '0.1514'
$ws.Range('E43').Value = '  -0.69%  '
$ws.Range('D44').Value = '111.55'
$ws.Range('E44').Value = '  +3.97%  '
$ws.Range('D45').Value = '10.59'
$ws.Range('E45').Value = '  -0.90%  '
$ws.Range('D46').Value = '0.4804'
$ws.Range('E46').Value = '  -0.61%  '
$ws.Range('D47').Value = '1.004'
$ws.Range('E47').Value = '  +0.12%  '
$ws.Range('D48').Value = '1.632'
$ws.Range('E48').Value = '  -1.57%  '
$ws.Range('D49').Value = '67.39'
$ws.Range('E49').Value = '  -1.37%  '
$ws.Range('D50').Value = '0.06084'
$ws.Range('E50').Value = '  -0.38%  '
$ws.Range('D51').Value = '0.9004'
$ws.Range('E51').Value = '  -0.08%  '

# Restore the default (Normal) style on the forced-text cells so no stray
# number-format style lingers on them, while keeping their content as text.
foreach ($cellRef in $textCells) {
    $ws.Range($cellRef).Style = "Normal"
}
